$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

# Row 2: now holds the identifiers that previously lived on row 3,
# and the seat changes from "C2" to "B1"
$ws.Range("A2").Value = "ff87f03b-8891-4bb6-ac5c-a510d216fdd6"
$ws.Range("C2").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$ws.Range("E2").Value = "B1"

# Row 3: now holds the identifiers that previously lived on row 2,
# and the seat changes from "B1" to "B2"
$ws.Range("A3").Value = "91e25164-6f67-42f7-b978-9132a406c060"
$ws.Range("C3").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$ws.Range("E3").Value = "B2"
